# Correction type pour génération à partir fsh
# Sets the "Name" property value on the Metadata sheet, and refreshes
# the "Date" property value to the new generation timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 4 -> A4 = "Name", B4 was blank; fill in the generated name.
$ws.Range("B4").Value = "TypesavoirfaireVs"

# Row 8 -> A8 = "Date", B8 holds the generation timestamp; bump it.
$ws.Range("B8").Value = "2025-07-18T06:40:38+00:00"
